# Normalization for CPF data: add Googleworkspace, Gympass, Unimed benefit
# columns (new F, G, H) and shift the Total column to I, recalculating totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the three new benefit headers, copying the existing header style
# (bold, bordered) from F1 onto the newly used G1:I1 range first.
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Googleworkspace"
$ws.Range("G1").Value = "Gympass"
$ws.Range("H1").Value = "Unimed"
$ws.Range("I1").Value = "Total"

# Per-employee data: row, Github(E, $null = unchanged), Googleworkspace(F),
# Gympass(G), Unimed(H), new Total(I)
$data = @(
    @(2, $null, 297.13, 90, 533.99, 4560.710765580439),
    @(3, $null, 297.13, 90, 855.8, 4962.033451680459),
    @(4, $null, 297.13, 90, 444.99, 9936.309721841499),
    @(5, $null, 297.13, 90, 444.99, 8891.909738823022),
    @(6, $null, 297.13, 90, 444.99, 8196.646801321949),
    @(7, $null, 297.13, 90, 444.99, 6151.107724937523),
    @(8, $null, 297.13, 90, 364.76, 10557.73905794158),
    @(9, $null, 297.13, 90, 560.7, 9486.596836853116),
    @(10, $null, 297.13, 90, 444.99, 2401.868209620668),
    @(11, $null, 297.13, 90, 533.99, 2399.06255272744),
    @(12, 0, 297.13, 90, 848.05, 5257.013629163991),
    @(13, 0, 297.13, 117, 533.99, 6192.743025188975),
    @(14, 0, 297.13, 90, 444.99, 8889.377421621037),
    @(15, 141.46, 297.13, 90, 533.99, 10775.30828538261),
    @(16, 141.46, 297.13, 90, 533.99, 8550.579892436921),
    @(17, $null, 297.13, 90, 444.99, 9003.11768092652),
    @(18, $null, 297.13, 90, 364.76, 2788.404997025348),
    @(19, $null, 297.13, 90, 444.99, 1940.460345861017),
    @(20, $null, 297.13, 90, 444.99, 11063.28102463056),
    @(21, $null, 297.13, 117, 560.7, 4524.492751492262),
    @(22, $null, 297.13, 90, 533.99, 9447.590794845279),
    @(23, $null, 297.13, 90, 444.99, 6120.380196934793),
    @(24, $null, 297.13, 90, 560.7, 6166.00282637831),
    @(25, $null, 297.13, 117, 364.76, 5366.657886848967),
    @(26, $null, 297.13, 90, 533.99, 4151.30657849905),
    @(27, $null, 297.13, 90, 444.99, 2194.351969577254),
    @(28, $null, 297.13, 206.57, 364.76, 8930.373739691342),
    @(29, $null, 297.13, 90, 364.76, 7139.183051889353),
    @(30, $null, 297.13, 90, 364.76, 7021.621383248207),
    @(31, $null, 297.13, 90, 533.99, 8442.306797367171),
    @(32, $null, 297.13, 90, 364.76, 5616.914047610113),
    @(33, $null, 297.13, 246, 364.76, 9350.35960081588),
    @(34, $null, 297.13, 90, 444.99, 2367.34645381395),
    @(35, $null, 297.13, 90, 364.76, 9844.295201541867),
    @(36, $null, 297.13, 164, 364.76, 4092.547357480083),
    @(37, 141.46, 297.13, 90, 364.76, 9391.195668097969),
    @(38, 0, 297.13, 90, 444.99, 4322.76446620121),
    @(39, 0, 297.13, 90, 364.76, 1964.057158944151),
    @(40, 254.63, 297.13, 90, 404.91, 1304.61215616471),
    @(41, 0, 297.13, 90, 0, 6289.564171124831),
    @(42, 254.63, 297.13, 732, 364.76, 11573.69931350882),
    @(43, 0, 297.13, 90, 444.99, 4829.201814681695),
    @(44, 0, 297.13, 90, 533.99, 6166.07395791246),
    @(45, 0, 297.13, 90, 364.76, 6935.66422031294),
    @(46, 0, 297.13, 90, 444.99, 7005.523303137569),
    @(47, 0, 297.13, 90, 684.63, 8136.911216079629),
    @(48, 0, 297.13, 90, 855.8, 10338.02015171554),
    @(49, 254.63, 297.13, 90, 560.7, 1949.411078367926),
    @(50, 254.63, 297.13, 90, 855.8, 10231.22129096568),
    @(51, 254.63, 297.13, 164, 364.76, 5032.724962932551),
    @(52, 0, 297.13, 90, -328.28, 7051.98850131002),
    @(53, 254.63, 297.13, 90, 622.41, 1754.344359132607),
    @(54, 0, 297.13, 90, 444.99, 2293.363910515265),
    @(55, 0, 297.13, 90, 444.99, 3573.485342728207),
    @(56, 0, 297.13, 90, 444.99, 5287.179040524838),
    @(57, 254.63, 297.13, 90, 364.76, 1806.700426414696),
    @(58, 0, 297.13, 90, 364.76, 8088.792571199591),
    @(59, 0, 297.13, 117, 560.7, 2703.459771639805),
    @(60, 254.63, 297.13, 90, 533.99, 7136.611253655476),
    @(61, 0, 297.13, 90, -17.8, 7675.728876912478),
    @(62, 254.63, 297.13, 90, 533.99, 5231.258897212339),
    @(63, $null, 297.13, 90, 444.99, 7855.791722399166)
)

foreach ($row in $data) {
    $r = $row[0]
    $github = $row[1]
    $google = $row[2]
    $gympass = $row[3]
    $unimed = $row[4]
    $total = $row[5]
    if ($null -ne $github) {
        $ws.Cells.Item($r, 5).Value = $github
    }
    $ws.Cells.Item($r, 6).Value = $google
    $ws.Cells.Item($r, 7).Value = $gympass
    $ws.Cells.Item($r, 8).Value = $unimed
    $ws.Cells.Item($r, 9).Value = $total
}
